# --------------------------------------------------------------------
# Edit 1: turn the first paragraph
#   "This is a Microsoft word document."
# into
#   "This is a Microsoft word document.  " (black, two trailing spaces)
#   + three red runs reconstructing:
#     "(This is a change – Version for main branch)"
# --------------------------------------------------------------------
$d = $word.ActiveDocument

$p1 = $d.Paragraphs(1).Range
$p1.InsertAfter("  ")

# Re-fetch paragraph range after the mutation so Start/End are current,
# then back up one character to land just before the paragraph mark.
$fresh = $d.Paragraphs(1).Range
$insertAt = $fresh.End - 1

$ins1 = $d.Range($insertAt, $insertAt)
$ins1.InsertAfter([string]::Concat("(This is a change ", [char]0x2013, " Ve"))
$ins1.Font.Color = 255

$ins2 = $d.Range($ins1.End, $ins1.End)
$ins2.InsertAfter("rsion for main branch")
$ins2.Font.Color = 255

$ins3 = $d.Range($ins2.End, $ins2.End)
$ins3.InsertAfter(")")
$ins3.Font.Color = 255

# --------------------------------------------------------------------
# Edit 2: delete the trailing paragraph
#   "ank God almighty, we are free at last."
# (the whole paragraph, including its mark, disappears; the preceding
#  paragraph "...Shall be lifted-nevermore!" becomes the last one and
#  keeps/absorbs the section properties).
# --------------------------------------------------------------------
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($count)
$lastPara.Range.Delete()

# --------------------------------------------------------------------
# Edit 3: now that no paragraph uses them any longer, the handful of
# custom / unused styles that Word prunes on save disappear from
# styles.xml. Locate them by name first (collection indices are only
# valid for the snapshot taken at scan time), then delete starting
# from the highest index so earlier lookups stay valid.
# --------------------------------------------------------------------
$staleStyleNames = @(
    "Heading 2",
    "Heading 4",
    "apple-converted-space",
    "Hyperlink",
    "Heading 2 Char",
    "Heading 4 Char",
    "audio-tool",
    "subscribe",
    "subscribe-more-info",
    "generic-title",
    "podcast-tools__subscribe-links"
)

$staleIndices = @()
for ($i = 1; $i -le $d.Styles.Count; $i++) {
    if ($staleStyleNames -contains $d.Styles.Item($i).NameLocal) {
        $staleIndices += $i
    }
}

foreach ($i in ($staleIndices | Sort-Object -Descending)) {
    $d.Styles.Item($i).Delete()
}
